# Add extension markup menu for genre and rework enum
#
# The "Films" sheet had every row's genre column (C) hard-coded to the
# single value "Історичне". This reworks it into a proper per-movie genre
# enum, replacing the old single shared string with a set of distinct
# genre values (one per row), widens column C so the longer genre labels
# fit, and leaves the selection on the genre column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Films")

# Assign the new genre values. The order below matches the order in which
# the distinct genre strings were first introduced, so the shared-strings
# table ends up built in the same sequence as the source workbook.
$ws.Cells.Item(4, 3).Value  = "Комедія"
$ws.Cells.Item(5, 3).Value  = "Документальний фільм"
$ws.Cells.Item(2, 3).Value  = "Містика"
$ws.Cells.Item(6, 3).Value  = "Фантастика"
$ws.Cells.Item(1, 3).Value  = "Історична драма"
$ws.Cells.Item(9, 3).Value  = "Анімація"
$ws.Cells.Item(10, 3).Value = "Екшн"
$ws.Cells.Item(8, 3).Value  = "Трилер"
$ws.Cells.Item(3, 3).Value  = "Романтика"
$ws.Cells.Item(7, 3).Value  = "Фантастика"

# Widen the genre column so the new (longer) labels are readable.
$ws.Columns.Item(3).ColumnWidth = 30.83

# Leave the selection on the genre column, as in the source edit.
$ws.Range("C12").Select()
